# Updated symbol list on Fri Dec 30 05:22:38 UTC 2022 with GitHub Actions
#
# Applies the latest coinranking.com price/volume snapshot to the
# "cryptos" worksheet. Most rows only get a refreshed Price (column D)
# value; a few rows also had their 24h best/worst badge text toggled in
# column E, and three rows (41-43) were re-sorted with fresh Coin/Link/
# Price/Volume data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Sheet, $Cell, $Text)
    # Force a plain-text number format so numeric-looking strings (e.g.
    # "245.57" or "0.006280") are kept verbatim as text rather than being
    # coerced into floating point numbers that would lose trailing zeros.
    $range = $Sheet.Range($Cell)
    $range.NumberFormat = "@"
    $range.Value = $Text
}

# --- Simple Price (column D) refreshes -----------------------------------
Set-TextValue $ws "D2"  "245.57"
Set-TextValue $ws "D4"  "5.244"
Set-TextValue $ws "D5"  "0.05794"
Set-TextValue $ws "D6"  "6.510"
Set-TextValue $ws "D7"  "3.123"
Set-TextValue $ws "D8"  "0.8156"
Set-TextValue $ws "D9"  "0.8519"
Set-TextValue $ws "D10" "0.1360"
Set-TextValue $ws "D11" "0.06960"
Set-TextValue $ws "D12" "0.03195"
Set-TextValue $ws "D13" "0.02872"
Set-TextValue $ws "D14" "0.09376"
Set-TextValue $ws "D15" "3.752"
Set-TextValue $ws "D16" "0.001517"
Set-TextValue $ws "D17" "0.04694"

# Row 18 also lost its "Worst in 24h" badge in column E.
Set-TextValue $ws "D18" "0.0005981"
Set-TextValue $ws "E18" "17OneONE"

Set-TextValue $ws "D19" "0.006283"
Set-TextValue $ws "D20" "0.001238"
Set-TextValue $ws "D21" "0.004537"
Set-TextValue $ws "D22" "0.00006903"
Set-TextValue $ws "D24" "2.085"
Set-TextValue $ws "D25" "0.3189"
Set-TextValue $ws "D28" "0.0002329"
Set-TextValue $ws "D40" "0.03662"

# --- Rows 41-43 were re-sorted with refreshed data ------------------------
# Old order: KickToken, BKEXToken, CEJI
# New order: BKEXToken, CEJI, KickToken
Set-TextValue $ws "B41" "BKEXToken"
Set-TextValue $ws "C41" "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue $ws "D41" "0.1054"
Set-TextValue $ws "E41" "40BKEXTokenBKK"

Set-TextValue $ws "B42" "CEJI"
Set-TextValue $ws "C42" "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue $ws "D42" "0.002751"
Set-TextValue $ws "E42" "41CEJICEJIBestin24h"

Set-TextValue $ws "B43" "KickToken"
Set-TextValue $ws "C43" "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue $ws "D43" "0.006282"
Set-TextValue $ws "E43" "42KickTokenKICK"

Set-TextValue $ws "D44" "0.007954"
Set-TextValue $ws "D45" "0.00005272"
Set-TextValue $ws "D47" "0.3300"

# Row 48 gained a "Worst in 24h" badge in column E.
Set-TextValue $ws "D48" "0.002340"
Set-TextValue $ws "E48" "47BOLOBOLOWorstin24h"
